# "all updates to include tri proximity tables"
#
# Means sheet: the proximity-band columns (Within 1/3/5/10 miles of HFC
# production facility -> D:G) could no longer be computed for rows 2-10,
# so they now surface as #NUM! errors.
#
# Standard Deviations sheet: the corresponding proximity-band columns
# E:G (3/5/10 miles) collapse to 0 for rows 2-10, matching column D
# (1 mile), which was already 0.

$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsMeans.Range("D2:G10").Value = "#NUM!"

$wsSd = $wb.Worksheets.Item("Standard Deviations")
$wsSd.Range("E2:G10").Value = 0
